$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column F content (rows 2-6) shifts one column to the right into
# column G (same text/style for every row), except the row 5 task text is
# updated from "If done, make plan for Friday" to "If done, make plan for
# Saturday".

# Copy the whole F2:F6 block (values + formatting) over to G2:G6.
$ws.Range("F2:F6").Copy($ws.Range("G2:G6"))
$excel.CutCopyMode = 0

# Update the task text for row 5 (now living in G5) for the new plan.
$ws.Range("G5").Value = "If done, make plan for Saturday"

# Remove the old column F cells (values + formatting) now that everything
# has moved into column G.
$ws.Range("F2:F6").Clear()

# Update the active selection to match the authored change.
$ws.Range("G6").Select()
